$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6111739799904399
$ws.Range("C2").Value = 0.186433905212855
$ws.Range("D2").Value = 0.05180973644350217
$ws.Range("E2").Value = 0.1194009475103286
$ws.Range("F2").Value = 1.068408293977036
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.3152865101613713
$ws.Range("L2").Value = 0.1920005324638652
$ws.Range("O2").Value = 3.84824323216958
$ws.Range("B3").Value = 0.568764282928953
$ws.Range("C3").Value = 0.1867681177691267
$ws.Range("D3").Value = 0.05010185833346981
$ws.Range("E3").Value = 0.1189125529880499
$ws.Range("F3").Value = 1.069010116001721
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.2791099968645767
$ws.Range("L3").Value = 0.1849179742812765
$ws.Range("O3").Value = 3.865075140854202
$ws.Range("B4").Value = 0.5429188198678787
$ws.Range("C4").Value = 0.1869959736744917
$ws.Range("D4").Value = 0.04904136516670121
$ws.Range("E4").Value = 0.11867376500755
$ws.Range("F4").Value = 1.069978404952778
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.2568984416939912
$ws.Range("L4").Value = 0.1806710816900079
$ws.Range("O4").Value = 3.877423588116102
$ws.Range("B5").Value = 0.5324360502573882
$ws.Range("C5").Value = 0.1870945466769953
$ws.Range("D5").Value = 0.04860624778882539
$ws.Range("E5").Value = 0.1185918329572537
$ws.Range("F5").Value = 1.07052356743867
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.2478477703463966
$ws.Range("L5").Value = 0.178966105287671
$ws.Range("O5").Value = 3.882961979804122
$ws.Range("B6").Value = 0.5306984002395154
$ws.Range("C6").Value = 0.1871112608505463
$ws.Range("D6").Value = 0.04853381887863151
$ws.Range("E6").Value = 0.1185791573128405
$ws.Range("F6").Value = 1.070623186081534
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.2463449706152829
$ws.Range("L6").Value = 0.1786845472899472
$ws.Range("O6").Value = 3.883912203790317
$ws.Range("B7").Value = 0.5427772444623145
$ws.Range("C7").Value = 0.1869972798714734
$ws.Range("D7").Value = 0.04903550896827369
$ws.Range("E7").Value = 0.1186725977665049
$ws.Range("F7").Value = 1.069985147510693
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.2567763776388148
$ws.Range("L7").Value = 0.1806479837753301
$ws.Range("O7").Value = 3.877496230868815
$ws.Range("B8").Value = 0.5965111530577758
$ws.Range("C8").Value = 0.1865444558894715
$ws.Range("D8").Value = 0.05122333189688533
$ws.Range("E8").Value = 0.1192198786790648
$ws.Range("F8").Value = 1.06849153236017
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.3028129194117071
$ws.Range("L8").Value = 0.1895373665808222
$ws.Range("O8").Value = 3.853628928863259
$ws.Range("B9").Value = 0.7034031555616878
$ws.Range("C9").Value = 0.1858351020292162
$ws.Range("D9").Value = 0.05541885827133086
$ws.Range("E9").Value = 0.1207773541943808
$ws.Range("F9").Value = 1.070314351611046
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.3930818299906491
$ws.Range("L9").Value = 0.2077757740710808
$ws.Range("O9").Value = 3.822806458634375
$ws.Range("B10").Value = 0.7828429776630799
$ws.Range("C10").Value = 0.1854214326336674
$ws.Range("D10").Value = 0.05844276359182743
$ws.Range("E10").Value = 0.122216548862486
$ws.Range("F10").Value = 1.074553548128662
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 0.4593818799431233
$ws.Range("L10").Value = 0.221666609997726
$ws.Range("O10").Value = 3.809913828208778
$ws.Range("B11").Value = 0.8191752283343874
$ws.Range("C11").Value = 0.185256293436133
$ws.Range("D11").Value = 0.05980555961085088
$ws.Range("E11").Value = 0.1229352700254971
$ws.Range("F11").Value = 1.077112569984678
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 0.4895362296731776
$ws.Range("L11").Value = 0.2280925395393893
$ws.Range("O11").Value = 3.80616878254142
$ws.Range("B12").Value = 0.8329607982145149
$ws.Range("C12").Value = 0.1851970491289094
$ws.Range("D12").Value = 0.06031975711236015
$ws.Range("E12").Value = 0.1232166276965927
$ws.Range("F12").Value = 1.078172310048998
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 0.5009536692645042
$ws.Range("L12").Value = 0.2305412087999059
$ws.Range("O12").Value = 3.805055604291283
$ws.Range("B13").Value = 0.8299906238260064
$ws.Range("C13").Value = 0.1852096624744277
$ws.Range("D13").Value = 0.06020909875410752
$ws.Range("E13").Value = 0.1231556236104616
$ws.Range("F13").Value = 1.077940042181837
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 0.4984947879209471
$ws.Range("L13").Value = 0.2300131642281116
$ws.Range("O13").Value = 3.805281780236015
$ws.Range("B14").Value = 0.8203088313989326
$ws.Range("C14").Value = 0.1852513535594795
$ws.Range("D14").Value = 0.0598479004392658
$ws.Range("E14").Value = 0.1229582333128079
$ws.Range("F14").Value = 1.077197937707055
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 0.4904755794527489
$ws.Range("L14").Value = 0.2282936865906748
$ws.Range("O14").Value = 3.806071087383941
$ws.Range("B15").Value = 0.8143819942825132
$ws.Range("C15").Value = 0.1852773183727052
$ws.Range("D15").Value = 0.05962641261107393
$ws.Range("E15").Value = 0.1228385229406612
$ws.Range("F15").Value = 1.07675518878122
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 0.4855633906394701
$ws.Range("L15").Value = 0.2272424486491929
$ws.Range("O15").Value = 3.806594283937699
$ws.Range("B16").Value = 0.7804724406327921
$ws.Range("C16").Value = 0.1854326865453118
$ws.Range("D16").Value = 0.05835344216615113
$ws.Range("E16").Value = 0.1221708659630458
$ws.Range("F16").Value = 1.074399000504968
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 0.4574110601098766
$ws.Range("L16").Value = 0.2212488062662459
$ws.Range("O16").Value = 3.810201246253939
$ws.Range("B17").Value = 0.7597194195235772
$ws.Range("C17").Value = 0.1855338848280894
$ws.Range("D17").Value = 0.05756922209840809
$ws.Range("E17").Value = 0.1217776689311911
$ws.Range("F17").Value = 1.073115076265722
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 0.4401386463704284
$ws.Range("L17").Value = 0.2175992386419523
$ws.Range("O17").Value = 3.812957061874414
$ws.Range("B18").Value = 0.7478011938877955
$ws.Range("C18").Value = 0.1855942614920068
$ws.Range("D18").Value = 0.05711695668486527
$ws.Range("E18").Value = 0.1215575394537645
$ws.Range("F18").Value = 1.072435950239708
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.4302034883361046
$ws.Range("L18").Value = 0.2155101714932925
$ws.Range("O18").Value = 3.814741657846383
$ws.Range("B19").Value = 0.7437690614924293
$ws.Range("C19").Value = 0.1856150774124075
$ws.Range("D19").Value = 0.05696362154500179
$ws.Range("E19").Value = 0.1214840429436954
$ws.Range("F19").Value = 1.072216203310731
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.4268395404744183
$ws.Range("L19").Value = 0.2148045803689485
$ws.Range("O19").Value = 3.815380155073058
$ws.Range("B20").Value = 0.7619267183088994
$ws.Range("C20").Value = 0.1855228876691655
$ws.Range("D20").Value = 0.05765282832116725
$ws.Range("E20").Value = 0.1218189017433353
$ws.Range("F20").Value = 1.073245609234945
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 0.4419773821012711
$ws.Range("L20").Value = 0.2179866998598072
$ws.Range("O20").Value = 3.81264304951651
$ws.Range("B21").Value = 0.8231518719024677
$ws.Range("C21").Value = 0.1852390187595674
$ws.Range("D21").Value = 0.05995404391231318
$ws.Range("E21").Value = 0.1230159622218068
$ws.Range("F21").Value = 1.077413450248926
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 0.4928310563099672
$ws.Range("L21").Value = 0.2287983238898619
$ws.Range("O21").Value = 3.805830970204283
$ws.Range("B22").Value = 0.8633250742963696
$ws.Range("C22").Value = 0.1850726621136758
$ws.Range("D22").Value = 0.06144714676299401
$ws.Range("E22").Value = 0.1238518914585285
$ws.Range("F22").Value = 1.080666016076236
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 0.5260587187786712
$ws.Range("L22").Value = 0.2359535483607829
$ws.Range("O22").Value = 3.803156663553608
$ws.Range("B23").Value = 0.8418695465844053
$ws.Range("C23").Value = 0.1851597034080896
$ws.Range("D23").Value = 0.06065125260926862
$ws.Range("E23").Value = 0.1234008419602404
$ws.Range("F23").Value = 1.078881684840312
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 0.5083254050418304
$ws.Range("L23").Value = 0.2321265300238622
$ws.Range("O23").Value = 3.804421273950851
$ws.Range("B24").Value = 0.7609287581026933
$ws.Range("C24").Value = 0.1855278526420392
$ws.Range("D24").Value = 0.0576150343279096
$ws.Range("E24").Value = 0.1218002419333288
$ws.Range("F24").Value = 1.073186411438613
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 0.4411461052687855
$ws.Range("L24").Value = 0.2178115002387102
$ws.Range("O24").Value = 3.812784390768655
$ws.Range("B25").Value = 0.6743254259203582
$ws.Range("C25").Value = 0.1860080237744235
$ws.Range("D25").Value = 0.05429408431602667
$ws.Range("E25").Value = 0.1203042032760528
$ws.Range("F25").Value = 1.069312189369512
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.3686642073799362
$ws.Range("L25").Value = 0.2027555266340073
$ws.Range("O25").Value = 3.829432809506471
